# Adds new brainstorming notes to the "Experiment Ideas" sheet and new
# literature/experiment notes to the "Dissertation Structure" sheet.
# Cells are written in the same order the text was originally typed so the
# shared-string table comes out in the same order too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Experiment Ideas"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Experiment Ideas")

$ws1.Range("C5").Value = "stock prices"
$ws1.Range("C6").Value = "volatility"
$ws1.Range("C7").Value = "inflation"
$ws1.Range("C8").Value = "unemployment"
$ws1.Range("C9").Value = "gdp"
$ws1.Range("C10").Value = "interest rates"

$ws1.Range("B12").Value = "different benchmark models"

$ws1.Range("B16").Value = "Different data frequencies"
$ws1.Range("C17").Value = "minute"
$ws1.Range("C18").Value = "hour"
$ws1.Range("C19").Value = "day"
$ws1.Range("C20").Value = "week"
$ws1.Range("C21").Value = "month"
$ws1.Range("C22").Value = "year"

$ws1.Range("C26").Value = "stable"
$ws1.Range("C25").Value = "recession"
$ws1.Range("B24").Value = "Different economic conditions"

$ws1.Range("B14").Value = "different configuration of the models"

$ws1.Range("B15").Select()

# ---------------------------------------------------------------------
# Sheet 2: "Dissertation Structure"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Dissertation Structure")

$ws2.Range("S15").Value = "1045 time series"
$ws2.Range("S16").Value = "length between 80 and 126"
$ws2.Range("S17").Value = "from every time-series 18 points were held out"
$ws2.Range("S18").Value = "one-step-ahead forecasting"
$ws2.Range("S19").Value = "time-series with different features:"

$ws2.Range("T20").Value = "seasonality vs non seasonality"
$ws2.Range("T21").Value = "trend (linear or exponential or none)"

$ws2.Range("S22").Value = "data pre-processing:"
$ws2.Range("T23").Value = "log transformation"
$ws2.Range("T24").Value = "deseasonalization"
$ws2.Range("T26").Value = "scaling"
$ws2.Range("U25").Value = "autocorrelation with lag 12 months, using Bartlett formula for confidence"

$ws2.Range("R27").Value = "results"
$ws2.Range("R42").Value = "comments on the results"

$ws2.Range("T29").Value = "average rank"
$ws2.Range("T30").Value = "SMAPE-TOT (symmetric mean absolute percentage error)"
$ws2.Range("T31").Value = "FRAC-BEST (fraction-best) - fraction of time series for which a specific model beats all other models"
$ws2.Range("U32").Value = "SMAPE is used for computing this measure"

$ws2.Range("S28").Value = "metrics"

$ws2.Range("S33").Value = "experiment"
$ws2.Range("T34").Value = "different preprocessing methods"
$ws2.Range("T35").Value = "data from different domains"
$ws2.Range("T36").Value = "data with different properties"
$ws2.Range("U37").Value = "trend, no trend, seasonality, etc…"

$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollColumn = 7
$ws2.Range("P36").Select()
